$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 514) { $lastRow = 514 }

$ws.Range("C2:C$lastRow").Value = 45182
